$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column H
$ws.Range("H1").Value = "Added_cf_num"
$ws.Range("H1").Font.Bold = $true

# Formulas for H3:H42 = G{row} - 27133
$ws.Range("H3").Formula = "=G3-27133"
$ws.Range("H4:H42").Formula = "=G4-27133"

# Autofit the new column
$ws.Columns.Item(8).AutoFit()

# Selection per diff
$ws.Range("A3:H6").Select()

# Page setup
$ws.PageSetup.Orientation = 1
